# Initial implementation of the automatic application updater
# Update the "go out" confirmation strings (rephrased) and the
# "Social networks" -> "Social Media" label, then add the new
# updater-related translation rows (20-22), and move the selection
# to reflect the last edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Idiomas")

# Row 5: confirmation-to-exit strings were reworded
$ws.Range("A5").Value = "¿Estás seguro de que quieres salir?"
$ws.Range("B5").Value = "Are you sure you wanna go out?"

# Row 18 col B: "Social networks" -> "Social Media"
$ws.Range("B18").Value = "Social Media"

# New rows for the auto-updater feature strings
$ws.Range("A20").Value = "Hay una nueva versión disponible"
$ws.Range("B20").Value = "A new version is available"

$ws.Range("A21").Value = "¿Quieres actualizar a la ultima versión?"
$ws.Range("B21").Value = "Do you wanna update to the latest version?"

$ws.Range("A22").Value = "Más tarde "
$ws.Range("B22").Value = "Later"

# Move the active selection to A23, matching the post-edit cursor position
$ws.Range("A23").Select()
